# [Fix] Modificacion de perfil
#
# Renames the "Test10" profile to "Tes10" on both the "Server IP" and
# "Test Server" sheets, and adds a new "Tes11" profile row (with its IP /
# T24 string) right below the existing rows, plus a styled blank row
# underneath it. Also switches the active/selected sheet from "Modulos"
# to "Test Server".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Server IP"
$ws2 = $wb.Worksheets.Item(2)   # "Test Server"

# --- Step 1: add the new "Tes11" row to "Server IP" first -------------
# (doing this before the rename keeps the shared-string ordering that
# Excel produced: Tes11, 10.169.1.7, Tes10, t24tes11)
$ws1.Range("A5").Value = "Tes11"
$ws1.Range("B5").Value = "10.169.1.7"
$ws1.Range("B5").VerticalAlignment = -4160   # xlVAlignTop
$ws1.Range("B5").HorizontalAlignment = -4152 # xlRight

# --- Step 2: rename "Test10" -> "Tes10" on both sheets -----------------
$ws1.Range("A2").Value = "Tes10"
$ws2.Range("A2").Value = "Tes10"

# --- Step 3: add the new "Tes11" row to "Test Server" -------------------
$ws2.Range("A5").Value = "Tes11"
$ws2.Range("B5").Value = "t24tes11"

# --- Step 4: styled blank row underneath the new data on both sheets ---
$ws1.Range("A6").HorizontalAlignment = -4131 # xlLeft
$ws1.Range("B6").VerticalAlignment = -4160   # xlVAlignTop
$ws1.Range("B6").HorizontalAlignment = -4152 # xlRight
$ws2.Range("A6").HorizontalAlignment = -4131 # xlLeft

# --- Step 5: update selections / active sheet ---------------------------
$ws1.Activate()
$ws1.Range("A5:B5").Select()

$ws2.Activate()
